$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.270535349845886
$ws.Range("B1").Value = 1.638614773750305
$ws.Range("C1").Value = 2.227327823638916
$ws.Range("D1").Value = 6.290139198303223
$ws.Range("E1").Value = 2.998435497283936
